$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: 2024 August
$ws.Range("B12").Value = 31
$ws.Range("C12").Value = 25

# Row 13: 2024 September
$ws.Range("B13").Value = 12
$ws.Range("C13").Value = 23
